$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14) on the "Repayment schedule"
# sheet. This shifts the old N/O/P columns ("Late"/blank-heading/"Outstanding")
# one position to the right (-> O/P/Q) and leaves the new column N blank.
$ws.Columns("N").Insert()

# The newly inserted column picks up the width of its left neighbour (column M,
# "In Advance") as a custom (non-bestFit) width (target stored width 10.7109375
# characters; the closest value this engine can round-trip is ~10.6667).
$ws.Columns("N").ColumnWidth = 9.85

# Make "Repayment schedule" the active sheet/tab and select cell S8 on it,
# matching the new sheetView/selection state.
$ws.Activate()
$ws.Range("S8").Select() | Out-Null
